$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "nama" (column C) values to be unique per row instead of the
# shared "Kantor PUSDIKLAT Kendalisada" string.
$ws.Range("C2").Value = "Kantor PUSDIKLAT Kendalisada 1"
$ws.Range("C3").Value = "Kantor PUSDIKLAT Kendalisada 2"
$ws.Range("C4").Value = "Kantor PUSDIKLAT Kendalisada 3"

# Fix "harga" (column N) values for rows 3 and 4 back to 0.
$ws.Range("N3").Value = 0
$ws.Range("N4").Value = 0

# Column C needs to widen slightly to fit the new, longer text
# (author's original best-fit width was 30.140625 characters).
$ws.Columns.Item(3).ColumnWidth = 29.3

# Update the active selection to match the author's last position.
$ws.Range("O3").Select()
